# Commit: "Changed UDC name from ExcelUdc to iPortal"
#
# The workbook's "Attributes" sheet lists various app/value pairs; the
# rows describing the "udc" app currently carry the value "EXCELUDC"
# (shared-string) in column C, rows 57-66. Rename that value to "iPortal".
#
# The author's Excel session also ended up scrolled/selected differently
# on that sheet (selection C58:C66, active cell C58, window scrolled so
# row 43 is the first visible row) - reproduce that view state too.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Attributes")
$ws.Activate()

# Rename the udc value from "EXCELUDC" to "iPortal" for every row that
# references it (C57:C66).
$ws.Range("C57:C66").Value = "iPortal"

# Match the author's final view/selection state on the sheet.
$excel.ActiveWindow.ScrollRow = 43
$excel.ActiveWindow.ScrollColumn = 1
[void]$ws.Range("C58:C66").Select()
